$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the existing header style (bold, centered, bordered) from H1
# so no new style entries are introduced, matching the source diff.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New data cells (row 2): I2 = 8, J2 = 8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
